# "Add Spain to the Website"
# - Fix the Barcelona filename so it has the .xml extension
# - Rename "Orange Sevilla" / "Orange Valladolid" to "Orange ES Sevilla" / "Orange ES Valladolid"
#   for naming consistency with the other Spanish (ES) entries
# - Flip the "In HTML" flag to "Yes" for the newly-published rows (Algeria + Spain block)
# - Leave the selection on the rows that were just edited (E15:E20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the missing ".xml" on the Barcelona filename
$ws.Range("D16").Value = "Orange ES Barcelona_100Gbps_ES.xml"

# Add the "ES" country prefix to the ISP names for Sevilla and Valladolid
$ws.Range("B19").Value = "Orange ES Sevilla"
$ws.Range("B20").Value = "Orange ES Valladolid"

# Mark these rows as published / present in the HTML output
$ws.Range("E14:E20").Value = "Yes"

# Match the author's final selection in the sheet
$ws.Range("E15:E20").Select() | Out-Null
